$d = $word.ActiveDocument

# The page used to end with a "Requisitos" section followed by a blank
# line, a "Ver no Jupiter ..." line and a "(c) 2020 ... Jekyll ..."
# site-footer line (leftover boilerplate scraped from the site build).
# This edit removes that trailing boilerplate, keeping the single blank
# paragraph (and the page-break paragraph) that originally followed it.

# Find the paragraph that closes the "Requisitos" section - deletion
# starts right after it.
$startPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*LOQ4236*") {
        $startPara = $p
    }
}

# Find the paragraph holding the Jekyll/Github-pages copyright footer -
# deletion ends at (and includes) this paragraph's own mark.
$endPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Contact: luizeleno@usp.br*") {
        $endPara = $p
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $delStart = $startPara.Next().Range.Start
    $delEnd = $endPara.Range.End
    $r = $d.Range($delStart, $delEnd)
    $r.Delete()
}
